# Add a "NEW MAXIMUM STOCK LEVEL" column (F) to the stock offtake report.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell. Row 1 already carries s="1" (bold + centered) with
# customFormat="1", so a freshly written cell in that row inherits the
# row's style automatically - no explicit Font/Alignment calls needed.
$ws.Range("F1").Value = "NEW MAXIMUM STOCK LEVEL"

# Target column width is 28 (character units) in the saved XML. This
# engine's <col> serializer adds a fixed 0.8333333333333334 padding on
# top of whatever ColumnWidth is assigned, so back that padding out here
# to land exactly on 28.
$ws.Columns.Item(6).ColumnWidth = 28 - 0.8333333333333334

# Move/select the new header cell, like the original author did after
# adding the column.
$ws.Range("F2").Select()
